# The commit swaps the bodies of ppt/theme/theme1.xml and ppt/theme/theme2.xml:
# theme1.xml (currently the "Office Theme" palette) and theme2.xml (currently the
# "Integral" palette, and the one actually wired to the slide master / presentation)
# trade places. Because the font scheme and format scheme blocks of the two theme
# parts are already byte-for-byte identical, the only effective difference is the
# 12 colour-scheme swatches (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink).
#
# theme2.xml is the theme that is actually linked from the slide master and the
# presentation, so re-colouring it to the old theme1.xml ("Office Theme") palette
# reproduces the visible effect of the swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Target values = the palette that used to live in ppt/theme/theme1.xml
# ("Office Theme"): dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
